$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C21").Value = 282
$ws.Range("D21").Value = 243
$ws.Range("E21").Value = 39
$ws.Range("F21").Value = 69.62750716332378
